$d = $word.ActiveDocument
$v = [char]11

# 1. Title
$d.Content.Find.Execute("Unraveling the Enigma of Dark Matter", $true, $false, $false, $false, $false, $true, 1, $false, "The Captivating Realm of Biology: Unveiling the Wonders of Life", 2) | Out-Null

# 2. Author name
$d.Content.Find.Execute("Emily Carter", $true, $false, $false, $false, $false, $true, 1, $false, "Sarah Johnson", 2) | Out-Null

# 3. Email
$d.Content.Find.Execute("emily.carter@astronomicalsociety.org", $true, $false, $false, $false, $false, $true, 1, $false, "sarahjohnson@highschool.edu", 2) | Out-Null

# 4. Body paragraph (paragraph 5) - full rewrite, preserving leading run formatting
$bodyPara = $d.Paragraphs.Item(5)
$bodyStart = $bodyPara.Range.Start
$bodyEnd = $bodyPara.Range.End - 1
$bodyRange = $d.Range($bodyStart, $bodyEnd)
$bodyRange.Text = "Biology, the study of life, unraveling the mysteries of living organisms, from the microscopic cells to the interconnected ecosystems that shape our planet. It is an engrossing field that intertwines chemistry, physics, and mathematics to unveil the symphony of biological processes that governs life on Earth.${v}${v}As we delve into the fascinating world of biology, we discover a symphony of life, orchestrated by molecules, cells, and organisms, each playing a unique role in the intricate web of life. We learn about the remarkable diversity of organisms, ranging from microscopic bacteria to towering trees, and the intricate mechanisms that enable them to survive and thrive in their respective environments.${v}${v}Furthermore, biology reveals the incredible interconnectedness of all living things, showcasing the intricate relationships between organisms and their ecosystems. We explore the food chains and webs that sustain life, the cycles that recycle nutrients and energy, and the complex interactions that maintain the delicate balance of nature.${v}${v}Body:${v}${v}Biology encompasses a broad scope of subfields, each with its unique focus and perspective. Molecular biology investigates the structure and function of molecules, the building blocks of life. The study of cells, known as cell biology, delves into the inner workings of these fundamental units of life, uncovering the secrets of cellular processes like metabolism, reproduction, and communication.${v}${v}Additionally, organismal biology explores the diversity of life forms, ranging from single-celled organisms to complex multicellular organisms, and the adaptations that enable them to survive in various environments. Evolutionary biology unravels the history of life on Earth, tracing the incredible journey of species over millions of years, highlighting the remarkable adaptations and the relentless drive to survive.${v}${v}The field of ecology investigates the relationships between living organisms and their environment, examining how populations interact, how ecosystems function, and how human activities impact these delicate balances. Finally, biomedical sciences focus on understanding the causes, prevention, and treatment of diseases, enabling the development of life-saving therapies and interventions."

# 5. Summary paragraph (paragraph 7) - full rewrite, preserving leading run formatting
$sumPara = $d.Paragraphs.Item(7)
$sumStart = $sumPara.Range.Start
$sumEnd = $sumPara.Range.End - 1
$sumRange = $d.Range($sumStart, $sumEnd)
$sumRange.Text = "Biology is a vast and captivating field of study that unravels the mysteries of life, encompassing the diversity, complexity, and interconnectedness of living organisms. It reveals the intricate molecular mechanisms that govern cellular processes, the remarkable adaptations that enable organisms to thrive in diverse environments, and the delicate balance of ecosystems that sustain life. From the fundamentals of molecules to the grandeur of ecosystems, biology inspires us with its boundless wonders and challenges us to unravel the secrets of life itself."

# 6. Add a new empty paragraph at the very end of the document body (after Summary paragraph)
$endRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Host "Done"
